# Populate the empty "results" table (rows 4-6, columns B:M) with the
# measured throughput values and the computed ratios, then center-align
# the raw-measurement block (B:I) the way the source data arrived.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (Average) ---------------------------------------------------
$ws.Range("B4").Value = 9827.463639666219
$ws.Range("C4").Value = 3817.4951773066332
$ws.Range("D4").Value = 73611.763668766012
$ws.Range("E4").Value = 30069.905299915819
$ws.Range("F4").Value = 324711.00390696019
$ws.Range("G4").Value = 114405.56928857259
$ws.Range("H4").Value = 425548.54651092086
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1.0396484374999999
$ws.Range("K4").Value = 1.0039843749999999
$ws.Range("L4").Value = 1.0003881454467769
$ws.Range("M4").Value = 1.0000390052795409

# --- Row 5 (Std. Dev.) ---------------------------------------------------
$ws.Range("B5").Value = 20321.686671218999
$ws.Range("C5").Value = 14653.429316186241
$ws.Range("D5").Value = 22454.243585952503
$ws.Range("E5").Value = 21727.084205178642
$ws.Range("F5").Value = 86919.085309273054
$ws.Range("G5").Value = 15208.53658317422
$ws.Range("H5").Value = 38837.619615386277
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1.00908203125
$ws.Range("K5").Value = 1.0009179687500001
$ws.Range("L5").Value = 1.000090599060059
$ws.Range("M5").Value = 1.000009155273438

# --- Row 6 (Bittorent) ---------------------------------------------------
$ws.Range("B6").Value = 111.3545
$ws.Range("C6").Value = 13.7204
$ws.Range("D6").Value = 871.77170000000001
$ws.Range("E6").Value = 192.30179999999999
$ws.Range("F6").Value = 8960.1813999999995
$ws.Range("G6").Value = 1448.7623000000001
$ws.Range("H6").Value = 154724.18340000001
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 1.056546
$ws.Range("K6").Value = 1.0095050000000001
$ws.Range("L6").Value = 1.002785
$ws.Range("M6").Value = 1.002524

# The raw per-file-size measurements (columns B:I) came in center-aligned;
# apply that now that the cells are populated.
$ws.Range("B4:I6").HorizontalAlignment = -4108

# Leave the selection where the upload left it.
$ws.Range("K7").Select()
